# Auto-generated edit script applying the Kraken_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
# Row 39
$ws.Cells.Item(39, 8).Value = 454.1  # H39
$ws.Cells.Item(39, 9).Value = 36.8  # I39
$ws.Cells.Item(39, 10).Value = 871.4  # J39
$ws.Cells.Item(39, 11).Value = 110.4  # K39
$ws.Cells.Item(39, 12).Value = 2614.2  # L39
$ws.Cells.Item(39, 13).Value = 185.6  # M39
$ws.Cells.Item(39, 14).Value = -3206.2  # N39

# Row 40
$ws.Cells.Item(40, 8).Value = 5556.6665  # H40
$ws.Cells.Item(40, 10).Value = 6061.1113  # J40
$ws.Cells.Item(40, 12).Value = 6061.1113  # L40
$ws.Cells.Item(40, 14).Value = -6411.1113  # N40

# Row 41
$ws.Cells.Item(41, 8).Value = 5999.5  # H41
$ws.Cells.Item(41, 9).Value = 1000  # I41
$ws.Cells.Item(41, 10).Value = 7666  # J41
$ws.Cells.Item(41, 11).Value = 1000  # K41
$ws.Cells.Item(41, 12).Value = 7666  # L41
$ws.Cells.Item(41, 13).Value = -560  # M41
$ws.Cells.Item(41, 14).Value = -8546  # N41

# Row 111
$ws.Cells.Item(111, 8).Value = 0  # H111
$ws.Cells.Item(111, 9).Value = 0  # I111
$ws.Cells.Item(111, 11).Value = 0  # K111
$ws.Cells.Item(111, 13).Value = ""  # M111 (cleared)

$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32, 8).Value = 1962.1052  # H32
$ws.Cells.Item(32, 9).Value = 2043.3334  # I32
$ws.Cells.Item(32, 10).Value = 500  # J32
$ws.Cells.Item(32, 11).Value = 2043.3334  # K32
$ws.Cells.Item(32, 12).Value = 500  # L32
$ws.Cells.Item(32, 13).Value = -1756.3334  # M32
$ws.Cells.Item(32, 14).Value = -1074  # N32

# Row 63
$ws.Cells.Item(63, 8).Value = 0  # H63
$ws.Cells.Item(63, 9).Value = 0  # I63
$ws.Cells.Item(63, 10).Value = 0  # J63
$ws.Cells.Item(63, 11).Value = 0  # K63
$ws.Cells.Item(63, 12).Value = ""  # L63 (cleared)
$ws.Cells.Item(63, 13).Value = ""  # M63 (cleared)
$ws.Cells.Item(63, 14).Value = 0  # N63

# Row 66
$ws.Cells.Item(66, 8).Value = 0  # H66
$ws.Cells.Item(66, 9).Value = 0  # I66
$ws.Cells.Item(66, 10).Value = 0  # J66
$ws.Cells.Item(66, 11).Value = 0  # K66
$ws.Cells.Item(66, 12).Value = ""  # L66 (cleared)
$ws.Cells.Item(66, 13).Value = ""  # M66 (cleared)
$ws.Cells.Item(66, 14).Value = 0  # N66

# Row 92
$ws.Cells.Item(92, 8).Value = 67500  # H92
$ws.Cells.Item(92, 10).Value = 67500  # J92
$ws.Cells.Item(92, 12).Value = 67500  # L92
$ws.Cells.Item(92, 14).Value = -72492  # N92

# Row 95
$ws.Cells.Item(95, 8).Value = 40208  # H95
$ws.Cells.Item(95, 10).Value = 40208  # J95
$ws.Cells.Item(95, 12).Value = 40208  # L95
$ws.Cells.Item(95, 14).Value = -45700  # N95

$ws = $wb.Worksheets.Item(3)
# Row 54
$ws.Cells.Item(54, 8).Value = 2450  # H54
$ws.Cells.Item(54, 9).Value = 2450  # I54
$ws.Cells.Item(54, 11).Value = 2450  # K54
$ws.Cells.Item(54, 13).Value = -1966  # M54

# Row 132
$ws.Cells.Item(132, 8).Value = 99995  # H132
$ws.Cells.Item(132, 10).Value = 99995  # J132
$ws.Cells.Item(132, 12).Value = 99995  # L132
$ws.Cells.Item(132, 14).Value = -110115  # N132

# Row 134
$ws.Cells.Item(134, 8).Value = 5733.5  # H134
$ws.Cells.Item(134, 9).Value = 4673.6  # I134
$ws.Cells.Item(134, 11).Value = 14020.8  # K134
$ws.Cells.Item(134, 13).Value = -11485.8  # M134

$ws = $wb.Worksheets.Item(4)
# Row 7
$ws.Cells.Item(7, 8).Value = 258.375  # H7
$ws.Cells.Item(7, 9).Value = 146  # I7
$ws.Cells.Item(7, 10).Value = 325.8  # J7
$ws.Cells.Item(7, 11).Value = 146  # K7
$ws.Cells.Item(7, 12).Value = 325.8  # L7
$ws.Cells.Item(7, 13).Value = -33  # M7
$ws.Cells.Item(7, 14).Value = -551.8  # N7

# Row 16
$ws.Cells.Item(16, 8).Value = 701.5  # H16
$ws.Cells.Item(16, 9).Value = 668.3333  # I16
$ws.Cells.Item(16, 10).Value = 1000  # J16
$ws.Cells.Item(16, 11).Value = 668.3333  # K16
$ws.Cells.Item(16, 12).Value = 1000  # L16
$ws.Cells.Item(16, 13).Value = -381.3333  # M16
$ws.Cells.Item(16, 14).Value = -1574  # N16

# Row 59
$ws.Cells.Item(59, 8).Value = 32000  # H59
$ws.Cells.Item(59, 9).Value = 20000  # I59
$ws.Cells.Item(59, 11).Value = 20000  # K59
$ws.Cells.Item(59, 13).Value = -18855  # M59

# Row 113
$ws.Cells.Item(113, 8).Value = 701.5  # H113
$ws.Cells.Item(113, 9).Value = 668.3333  # I113
$ws.Cells.Item(113, 10).Value = 1000  # J113
$ws.Cells.Item(113, 11).Value = 668.3333  # K113
$ws.Cells.Item(113, 12).Value = 1000  # L113
$ws.Cells.Item(113, 13).Value = 1501.6667  # M113
$ws.Cells.Item(113, 14).Value = -5340  # N113

$ws = $wb.Worksheets.Item(5)
# Row 39
$ws.Cells.Item(39, 8).Value = 1914  # H39
$ws.Cells.Item(39, 9).Value = 2155.5  # I39
$ws.Cells.Item(39, 10).Value = 1833.5  # J39
$ws.Cells.Item(39, 11).Value = 6466.5  # K39
$ws.Cells.Item(39, 12).Value = 5500.5  # L39
$ws.Cells.Item(39, 13).Value = -6172.5  # M39
$ws.Cells.Item(39, 14).Value = -6088.5  # N39

# Row 55
$ws.Cells.Item(55, 8).Value = 1664  # H55
$ws.Cells.Item(55, 10).Value = 4004.5  # J55
$ws.Cells.Item(55, 12).Value = 12013.5  # L55
$ws.Cells.Item(55, 14).Value = -12367.5  # N55

# Row 69
$ws.Cells.Item(69, 8).Value = 500  # H69
$ws.Cells.Item(69, 9).Value = 500  # I69
$ws.Cells.Item(69, 11).Value = 1500  # K69
$ws.Cells.Item(69, 13).Value = -689  # M69

# Row 72
$ws.Cells.Item(72, 8).Value = 500  # H72
$ws.Cells.Item(72, 9).Value = 500  # I72
$ws.Cells.Item(72, 11).Value = 4500  # K72
$ws.Cells.Item(72, 13).Value = -444  # M72

# Row 75
$ws.Cells.Item(75, 8).Value = 733  # H75
$ws.Cells.Item(75, 10).Value = 733  # J75
$ws.Cells.Item(75, 12).Value = 2199  # L75
$ws.Cells.Item(75, 14).Value = -4195  # N75

# Row 78
$ws.Cells.Item(78, 8).Value = 733  # H78
$ws.Cells.Item(78, 10).Value = 733  # J78
$ws.Cells.Item(78, 12).Value = 6597  # L78
$ws.Cells.Item(78, 14).Value = -16581  # N78

# Row 103
$ws.Cells.Item(103, 8).Value = 2005  # H103
$ws.Cells.Item(103, 9).Value = 923.75  # I103
$ws.Cells.Item(103, 10).Value = 2725.8333  # J103
$ws.Cells.Item(103, 11).Value = 2771.25  # K103
$ws.Cells.Item(103, 12).Value = 8177.499899999999  # L103
$ws.Cells.Item(103, 13).Value = -1892.25  # M103
$ws.Cells.Item(103, 14).Value = -9935.499899999999  # N103

# Row 136
$ws.Cells.Item(136, 8).Value = 3921.875  # H136
$ws.Cells.Item(136, 9).Value = 3921.875  # I136
$ws.Cells.Item(136, 11).Value = 11765.625  # K136
$ws.Cells.Item(136, 13).Value = -6665.625  # M136

$ws = $wb.Worksheets.Item(6)
# Row 2
$ws.Cells.Item(2, 8).Value = 288.4  # H2
$ws.Cells.Item(2, 9).Value = 273  # I2
$ws.Cells.Item(2, 11).Value = 273  # K2
$ws.Cells.Item(2, 13).Value = -160  # M2

# Row 92
$ws.Cells.Item(92, 8).Value = 7105  # H92
$ws.Cells.Item(92, 10).Value = 7105  # J92
$ws.Cells.Item(92, 12).Value = 7105  # L92
$ws.Cells.Item(92, 14).Value = -10849  # N92

# Row 113
$ws.Cells.Item(113, 8).Value = 1100  # H113
$ws.Cells.Item(113, 9).Value = 1100  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 1100  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).Value = ""  # M113 (cleared)
$ws.Cells.Item(113, 14).Value = 1070  # N113

# Row 132
$ws.Cells.Item(132, 8).Value = 3352.9092  # H132
$ws.Cells.Item(132, 9).Value = 3352.9092  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 10058.7276  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = ""  # M132 (cleared)
$ws.Cells.Item(132, 14).Value = -7528.7276  # N132

$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 8000  # H7
$ws.Cells.Item(7, 9).Value = 8000  # I7
$ws.Cells.Item(7, 10).Value = 0  # J7
$ws.Cells.Item(7, 11).Value = 8000  # K7
$ws.Cells.Item(7, 12).Value = 0  # L7
$ws.Cells.Item(7, 13).Value = ""  # M7 (cleared)
$ws.Cells.Item(7, 14).Value = -7888  # N7

# Row 16
$ws.Cells.Item(16, 8).Value = 4321  # H16
$ws.Cells.Item(16, 9).Value = 4321  # I16
$ws.Cells.Item(16, 10).Value = 0  # J16
$ws.Cells.Item(16, 11).Value = 4321  # K16
$ws.Cells.Item(16, 12).Value = 0  # L16
$ws.Cells.Item(16, 13).Value = ""  # M16 (cleared)
$ws.Cells.Item(16, 14).Value = -4151  # N16

# Row 46
$ws.Cells.Item(46, 8).Value = 3561.6  # H46
$ws.Cells.Item(46, 10).Value = 3735.111  # J46
$ws.Cells.Item(46, 12).Value = 3735.111  # L46
$ws.Cells.Item(46, 14).Value = -4111.111  # N46

# Row 55
$ws.Cells.Item(55, 8).Value = 1294.8334  # H55
$ws.Cells.Item(55, 9).Value = 929.875  # I55
$ws.Cells.Item(55, 11).Value = 929.875  # K55
$ws.Cells.Item(55, 13).Value = -756.875  # M55

# Row 121
$ws.Cells.Item(121, 8).Value = 73000  # H121
$ws.Cells.Item(121, 9).Value = 0  # I121
$ws.Cells.Item(121, 10).Value = 73000  # J121
$ws.Cells.Item(121, 11).Value = 0  # K121
$ws.Cells.Item(121, 12).Value = ""  # L121 (cleared)
$ws.Cells.Item(121, 13).Value = 73000  # M121
$ws.Cells.Item(121, 14).Value = -76494  # N121

# Row 122
$ws.Cells.Item(122, 8).Value = 8133.3335  # H122
$ws.Cells.Item(122, 10).Value = 8000  # J122
$ws.Cells.Item(122, 12).Value = 24000  # L122
$ws.Cells.Item(122, 14).Value = -28900  # N122

# Row 126
$ws.Cells.Item(126, 8).Value = 8000  # H126
$ws.Cells.Item(126, 9).Value = 8000  # I126
$ws.Cells.Item(126, 10).Value = 0  # J126
$ws.Cells.Item(126, 11).Value = 24000  # K126
$ws.Cells.Item(126, 12).Value = 0  # L126
$ws.Cells.Item(126, 13).Value = ""  # M126 (cleared)
$ws.Cells.Item(126, 14).Value = -21530  # N126

# Row 132
$ws.Cells.Item(132, 8).Value = 3355.4  # H132
$ws.Cells.Item(132, 10).Value = 1894  # J132
$ws.Cells.Item(132, 12).Value = 5682  # L132
$ws.Cells.Item(132, 14).Value = -10742  # N132

$ws = $wb.Worksheets.Item(8)
# Row 40
$ws.Cells.Item(40, 8).Value = 20000  # H40
$ws.Cells.Item(40, 9).Value = 20000  # I40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 11).Value = 20000  # K40
$ws.Cells.Item(40, 12).Value = 0  # L40
$ws.Cells.Item(40, 13).Value = ""  # M40 (cleared)
$ws.Cells.Item(40, 14).Value = -19851  # N40

# Row 100
$ws.Cells.Item(100, 8).Value = 15959.6  # H100
$ws.Cells.Item(100, 9).Value = 15959.6  # I100
$ws.Cells.Item(100, 11).Value = 31919.2  # K100
$ws.Cells.Item(100, 13).Value = -31378.2  # M100

# Row 113
$ws.Cells.Item(113, 8).Value = 1102.2858  # H113
$ws.Cells.Item(113, 9).Value = 415.25  # I113
$ws.Cells.Item(113, 10).Value = 2018.3334  # J113
$ws.Cells.Item(113, 11).Value = 1245.75  # K113
$ws.Cells.Item(113, 12).Value = 6055.0002  # L113
$ws.Cells.Item(113, 13).Value = 924.25  # M113
$ws.Cells.Item(113, 14).Value = -10395.0002  # N113

# Row 132
$ws.Cells.Item(132, 8).Value = 10070.375  # H132
$ws.Cells.Item(132, 9).Value = 8610.5  # I132
$ws.Cells.Item(132, 10).Value = 14450  # J132
$ws.Cells.Item(132, 11).Value = 25831.5  # K132
$ws.Cells.Item(132, 12).Value = 43350  # L132
$ws.Cells.Item(132, 13).Value = -23301.5  # M132
$ws.Cells.Item(132, 14).Value = -48410  # N132

